# Rename the diff-table column headers from *_old/_new to the
# format-version-specific suffixes *_FV2404/_FV2410, then turn the
# header range into a proper Excel Table (with autofilter) and freeze
# the header row - matching the upstream commit
# "chore: adapt column header formatting to respective input file names".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J: "<Name>_old" -> "<Name>_FV2404"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2404"
}

# Column K: "diff" stays as-is
$ws.Cells.Item(1, 11).Value = "diff"

# Columns L-U: "<Name>_new" -> "<Name>_FV2410"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2410"
}

# Turn A1:U94 into a native Excel Table (adds xl/tables/table1.xml,
# the autoFilter, and the <tableParts> reference on the worksheet).
$dataRange = $ws.Range("A1:U94")
$table = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$table.Name = "Table1"

# Freeze the header row (row 1) so it stays visible while scrolling.
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
